$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual F-column (imputed/missing) values in rows 2-25 ---
$ws.Cells.Item(6, 6).Value = 16.43      # F6: was missing -> 16.43
$ws.Cells.Item(8, 6).ClearContents()    # F8: 17.05 -> missing
$ws.Cells.Item(18, 6).Value = 18.35     # F18: was missing -> 18.35
$ws.Cells.Item(20, 6).ClearContents()   # F20: 17.73 -> missing
$ws.Cells.Item(23, 6).Value = 16.48     # F23: was missing -> 16.48
$ws.Cells.Item(25, 6).ClearContents()   # F25: 16.6 -> missing

# --- Remove the "RM 232" row entirely (row 26) ---
$ws.Rows.Item(26).Delete()

# --- Remove what is now the "SC 92" row (originally row 28, now row 27) ---
$ws.Rows.Item(27).Delete()

# --- Fix up the A-column values (column header "A") that shifted with the
#     rows above and whose "missing" pattern changed in the re-sampled data ---
$ws.Cells.Item(27, 2).Value = -20.4     # SC 101 (row27) B: missing -> -20.4
$ws.Cells.Item(28, 2).ClearContents()   # SC 105 (row28) B: -19.6 -> missing
$ws.Cells.Item(29, 2).ClearContents()   # SC 119 (row29) B: -19.5 -> missing
$ws.Cells.Item(30, 2).Value = -19.7     # SC 120 (row30) B: missing -> -19.7
$ws.Cells.Item(30, 6).Value = 16.89     # SC 120 (row30) F: missing -> 16.89
$ws.Cells.Item(32, 2).ClearContents()   # SC 193 (row32) B: -19.9 -> missing
